$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF" - same style as other headers (H1 uses style index 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for I2:J12
$data = @(
    @(10, 11),
    @(9, 9),
    @(7, 8),
    @(5, 7),
    @(8, 9),
    @(9, 9),
    @(6, 6),
    @(5, 6),
    @(6, 7),
    @(4, 6),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
